$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'FAPs'
$ws.Cells.Item(2, 2).Value = 'Wnt5b'
$ws.Cells.Item(2, 3).Value = 'Fzd2'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 7).Value = [double]"0.3693360000000001"
$ws.Cells.Item(2, 8).Value = [double]"1.108008"
$ws.Cells.Item(2, 9).Value = [double]"0.9453581798061689"
$ws.Cells.Item(2, 10).Value = [double]"0.9453581798061688"
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = [double]"0.3007906666666667"
$ws.Cells.Item(2, 14).Value = [double]"0.902372"
$ws.Cells.Item(2, 15).Value = [double]"0.03537029821880876"
$ws.Cells.Item(2, 16).Value = [double]"0.03537029821880876"
$ws.Cells.Item(2, 17).Value = [double]"0.111092821664"
$ws.Cells.Item(2, 18).Value = [double]"0.999835394976"
$ws.Cells.Item(2, 19).Value = [double]"0.03343760074333443"
$ws.Cells.Item(2, 20).Value = [double]"0.03343760074333443"

# Row 3
$ws.Cells.Item(3, 1).Value = 'FAPs'
$ws.Cells.Item(3, 2).Value = 'Wnt5b'
$ws.Cells.Item(3, 3).Value = 'Fzd2'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3, 7).Value = [double]"0.3693360000000001"
$ws.Cells.Item(3, 8).Value = [double]"1.108008"
$ws.Cells.Item(3, 9).Value = [double]"0.9453581798061689"
$ws.Cells.Item(3, 10).Value = [double]"0.9453581798061688"
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = [double]"7.647982"
$ws.Cells.Item(3, 14).Value = [double]"22.943946"
$ws.Cells.Item(3, 15).Value = [double]"0.899334434508434"
$ws.Cells.Item(3, 16).Value = [double]"0.899334434508434"
$ws.Cells.Item(3, 17).Value = [double]"2.824675079952"
$ws.Cells.Item(3, 18).Value = [double]"25.422075719568"
$ws.Cells.Item(3, 19).Value = [double]"0.8501931640439033"
$ws.Cells.Item(3, 20).Value = [double]"0.8501931640439032"

# Row 4
$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 2).Value = 'Wnt5b'
$ws.Cells.Item(4, 3).Value = 'Fzd2'
$ws.Cells.Item(4, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 7).Value = [double]"0.3693360000000001"
$ws.Cells.Item(4, 8).Value = [double]"1.108008"
$ws.Cells.Item(4, 9).Value = [double]"0.9453581798061689"
$ws.Cells.Item(4, 10).Value = [double]"0.9453581798061688"
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 13).Value = [double]"0.0008990000000000001"
$ws.Cells.Item(4, 14).Value = [double]"0.002697"
$ws.Cells.Item(4, 15).Value = [double]"0.0001057143775473167"
$ws.Cells.Item(4, 16).Value = [double]"0.0001057143775473167"
$ws.Cells.Item(4, 17).Value = [double]"0.0003320330640000001"
$ws.Cells.Item(4, 18).Value = [double]"0.002988297576"
$ws.Cells.Item(4, 19).Value = [double]"9.993795153747344e-05"
$ws.Cells.Item(4, 20).Value = [double]"9.993795153747342e-05"

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Wnt5b'
$ws.Cells.Item(5, 3).Value = 'Fzd2'
$ws.Cells.Item(5, 4).Value = 'MuSCs'
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 7).Value = [double]"0.3693360000000001"
$ws.Cells.Item(5, 8).Value = [double]"1.108008"
$ws.Cells.Item(5, 9).Value = [double]"0.9453581798061689"
$ws.Cells.Item(5, 10).Value = [double]"0.9453581798061688"
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = [double]"0.5528646666666667"
$ws.Cells.Item(5, 14).Value = [double]"1.658594"
$ws.Cells.Item(5, 15).Value = [double]"0.06501195117304938"
$ws.Cells.Item(5, 16).Value = [double]"0.06501195117304936"
$ws.Cells.Item(5, 17).Value = [double]"0.204192824528"
$ws.Cells.Item(5, 18).Value = [double]"1.837735420752"
$ws.Cells.Item(5, 19).Value = [double]"0.06145957982660149"
$ws.Cells.Item(5, 20).Value = [double]"0.06145957982660147"

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Wnt5b'
$ws.Cells.Item(6, 3).Value = 'Fzd2'
$ws.Cells.Item(6, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6, 7).Value = [double]"0.3693360000000001"
$ws.Cells.Item(6, 8).Value = [double]"1.108008"
$ws.Cells.Item(6, 9).Value = [double]"0.9453581798061689"
$ws.Cells.Item(6, 10).Value = [double]"0.9453581798061688"
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6, 13).Value = [double]"0.001510333333333333"
$ws.Cells.Item(6, 14).Value = [double]"0.004531"
$ws.Cells.Item(6, 15).Value = [double]"0.0001776017221605087"
$ws.Cells.Item(6, 16).Value = [double]"0.0001776017221605087"
$ws.Cells.Item(6, 17).Value = [double]"0.0005578204720000001"
$ws.Cells.Item(6, 18).Value = [double]"0.005020384248000001"
$ws.Cells.Item(6, 19).Value = [double]"0.0001678972407920994"
$ws.Cells.Item(6, 20).Value = [double]"0.0001678972407920994"

# Row 7
$ws.Cells.Item(7, 1).Value = 'MuSCs'
$ws.Cells.Item(7, 2).Value = 'Wnt5b'
$ws.Cells.Item(7, 3).Value = 'Fzd2'
$ws.Cells.Item(7, 4).Value = 'ECs'
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(7, 7).Value = [double]"0.02134766666666667"
$ws.Cells.Item(7, 8).Value = [double]"0.064043"
$ws.Cells.Item(7, 9).Value = [double]"0.05464182019383115"
$ws.Cells.Item(7, 10).Value = [double]"0.05464182019383115"
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = [double]"0.3007906666666667"
$ws.Cells.Item(7, 14).Value = [double]"0.902372"
$ws.Cells.Item(7, 15).Value = [double]"0.03537029821880876"
$ws.Cells.Item(7, 16).Value = [double]"0.03537029821880876"
$ws.Cells.Item(7, 17).Value = [double]"0.006421178888444444"
$ws.Cells.Item(7, 18).Value = [double]"0.057790609996"
$ws.Cells.Item(7, 19).Value = [double]"0.001932697475474335"
$ws.Cells.Item(7, 20).Value = [double]"0.001932697475474335"

# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Wnt5b'
$ws.Cells.Item(8, 3).Value = 'Fzd2'
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 7).Value = [double]"0.02134766666666667"
$ws.Cells.Item(8, 8).Value = [double]"0.064043"
$ws.Cells.Item(8, 9).Value = [double]"0.05464182019383115"
$ws.Cells.Item(8, 10).Value = [double]"0.05464182019383115"
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = [double]"7.647982"
$ws.Cells.Item(8, 14).Value = [double]"22.943946"
$ws.Cells.Item(8, 15).Value = [double]"0.899334434508434"
$ws.Cells.Item(8, 16).Value = [double]"0.899334434508434"
$ws.Cells.Item(8, 17).Value = [double]"0.1632665704086667"
$ws.Cells.Item(8, 18).Value = [double]"1.469399133678"
$ws.Cells.Item(8, 19).Value = [double]"0.04914127046453067"
$ws.Cells.Item(8, 20).Value = [double]"0.04914127046453067"

# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Wnt5b'
$ws.Cells.Item(9, 3).Value = 'Fzd2'
$ws.Cells.Item(9, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 7).Value = [double]"0.02134766666666667"
$ws.Cells.Item(9, 8).Value = [double]"0.064043"
$ws.Cells.Item(9, 9).Value = [double]"0.05464182019383115"
$ws.Cells.Item(9, 10).Value = [double]"0.05464182019383115"
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 13).Value = [double]"0.0008990000000000001"
$ws.Cells.Item(9, 14).Value = [double]"0.002697"
$ws.Cells.Item(9, 15).Value = [double]"0.0001057143775473167"
$ws.Cells.Item(9, 16).Value = [double]"0.0001057143775473167"
$ws.Cells.Item(9, 17).Value = [double]"1.919155233333334e-05"
$ws.Cells.Item(9, 18).Value = [double]"0.000172723971"
$ws.Cells.Item(9, 19).Value = [double]"5.77642600984326e-06"
$ws.Cells.Item(9, 20).Value = [double]"5.77642600984326e-06"

# Row 10
$ws.Cells.Item(10, 1).Value = 'MuSCs'
$ws.Cells.Item(10, 2).Value = 'Wnt5b'
$ws.Cells.Item(10, 3).Value = 'Fzd2'
$ws.Cells.Item(10, 4).Value = 'MuSCs'
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(10, 7).Value = [double]"0.02134766666666667"
$ws.Cells.Item(10, 8).Value = [double]"0.064043"
$ws.Cells.Item(10, 9).Value = [double]"0.05464182019383115"
$ws.Cells.Item(10, 10).Value = [double]"0.05464182019383115"
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = [double]"0.5528646666666667"
$ws.Cells.Item(10, 14).Value = [double]"1.658594"
$ws.Cells.Item(10, 15).Value = [double]"0.06501195117304938"
$ws.Cells.Item(10, 16).Value = [double]"0.06501195117304936"
$ws.Cells.Item(10, 17).Value = [double]"0.01180237061577778"
$ws.Cells.Item(10, 18).Value = [double]"0.106221335542"
$ws.Cells.Item(10, 19).Value = [double]"0.003552371346447895"
$ws.Cells.Item(10, 20).Value = [double]"0.003552371346447894"

# Row 11
$ws.Cells.Item(11, 1).Value = 'MuSCs'
$ws.Cells.Item(11, 2).Value = 'Wnt5b'
$ws.Cells.Item(11, 3).Value = 'Fzd2'
$ws.Cells.Item(11, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(11, 7).Value = [double]"0.02134766666666667"
$ws.Cells.Item(11, 8).Value = [double]"0.064043"
$ws.Cells.Item(11, 9).Value = [double]"0.05464182019383115"
$ws.Cells.Item(11, 10).Value = [double]"0.05464182019383115"
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(11, 13).Value = [double]"0.001510333333333333"
$ws.Cells.Item(11, 14).Value = [double]"0.004531"
$ws.Cells.Item(11, 15).Value = [double]"0.0001776017221605087"
$ws.Cells.Item(11, 16).Value = [double]"0.0001776017221605087"
$ws.Cells.Item(11, 17).Value = [double]"3.224209255555556e-05"
$ws.Cells.Item(11, 18).Value = [double]"0.000290178833"
$ws.Cells.Item(11, 19).Value = [double]"9.704481368409272e-06"
$ws.Cells.Item(11, 20).Value = [double]"9.704481368409274e-06"

Write-Host "edit complete; dimension now $($ws.UsedRange.Address())"
